# Updates cryptos price/volume data to match the latest scrape.
# Numeric-looking text (e.g. "136.85") must stay text, matching the source
# sheet where Price/Volume(1h) columns are plain strings (not numbers) --
# otherwise Excel auto-converts them to floats (losing formatting / precision).
# We force text via NumberFormat "@" then restore the default "Normal" style
# so the cell keeps looking like every other untouched cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.434.00"
$ws.Range("E2").Value = "  +0.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.109.40"
$ws.Range("E3").Value = "  +0.65%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.74%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.106.35"
$ws.Range("E8").Value = "  +0.63%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.447"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.27%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +1.21%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.36%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.66%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.648.44"
$ws.Range("E13").Value = "  +0.87%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +3.18%  "

# Row 15 - Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.38%  "

# Row 16 - ShibaInu
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000162"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "57.523.65"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.112.73"
$ws.Range("E18").Value = "  +0.69%  "

# Row 19 - Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.80%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.86%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.93%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "347.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.60%  "

# Row 23 - LEO
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

# Row 24 - Dai
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.90%  "

# Row 26 - Polygon
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.501"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.82%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -0.51%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  -0.19%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0902"
$ws.Range("E29").Value = "  -0.61%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.38%  "

# Row 31 - USDe
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.61%  "

# Row 33 - RenderToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.33%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "

# Row 35 - NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.01%  "

# Row 36 - Fetch.AI
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.38%  "

# Row 37 - Monero
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "

# Row 38 - Aptos
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.31%  "

# Row 39 - EnergySwap
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.75%  "

# Row 40 - ImmutableX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.31%  "

# Row 41 - Filecoin -> Stacks (row swap)
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.83%  "

# Row 42 - Stacks -> Filecoin (row swap)
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.21%  "

# Row 43 - Hedera
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0662"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "

# Row 44 - Mantle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.700"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.51%  "

# Row 45 - RenzoRestakedETH
$ws.Range("D45").Value = "3.148.27"
$ws.Range("E45").Value = "  +0.60%  "

# Row 46 - Maker -> OKB (row swap)
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "

# Row 47 - OKB -> Maker (row swap)
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.355.49"
$ws.Range("E47").Value = "  +2.40%  "

# Row 48 - FirstDigitalUSD
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "

# Row 49 - VeChain
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0267"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.30%  "

# Row 50 - ONDO
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.960"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "

# Row 51 - Cosmos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "

